$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 95

# Columns A and D hold values that Excel would otherwise auto-convert
# (A looks like a date, D looks like a number) -- force them to text so
# they are stored the same way as the existing rows (plain text values).
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2023-06-30"

$ws.Range("B$row").Value = "18:45:49"
$ws.Range("C$row").Value = "Friday"

$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "26"

$ws.Range("E$row").Value = 123463
$ws.Range("F$row").Value = 134483
$ws.Range("G$row").Value = 163727
$ws.Range("H$row").Value = 133884
$ws.Range("I$row").Value = 176811
$ws.Range("J$row").Value = 115748
$ws.Range("K$row").Value = 204946
$ws.Range("L$row").Value = 226104
$ws.Range("M$row").Value = 176538
$ws.Range("N$row").Value = 104615
$ws.Range("O$row").Value = 39818
$ws.Range("P$row").Value = 33688
$ws.Range("Q$row").Value = 52640
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 36156
$ws.Range("T$row").Value = -1
